$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new trading-day rows (2020-02-27 and 2020-02-28) to the
# historical price table, extending the used range from A1:I26 to A1:I28.
#
# Columns B (date) and C (id) hold text that LOOKS like a number/date
# ("2020-02-27", "0217"). A plain Range.Value assignment would make Excel
# auto-convert those into a date serial / numeric value, so each is
# briefly switched to a Text number format before the assignment, then
# restored to the workbook's default "Normal" style so the cell doesn't
# keep a stray number-format override (matching the rest of the sheet,
# where these text columns carry no explicit style).

# Row 27: 2020-02-27
$ws.Range("A27").Value = 1582761600
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "2020-02-27"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0217"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "PWRWELL"
$ws.Range("E27").Value = 0.38
$ws.Range("F27").Value = 0.395
$ws.Range("G27").Value = 0.35
$ws.Range("H27").Value = 0.355
$ws.Range("I27").Value = 97200800

# Row 28: 2020-02-28
$ws.Range("A28").Value = 1582848000
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "2020-02-28"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0217"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "PWRWELL"
$ws.Range("E28").Value = 0.345
$ws.Range("F28").Value = 0.35
$ws.Range("G28").Value = 0.295
$ws.Range("H28").Value = 0.32
$ws.Range("I28").Value = 63846900
